# Update the Spanish -> Portuguese market/trader name translations on the
# "main" sheet, refresh the C2:C15 selection, size column C to fit its new
# (longer) contents, and restore the default "Office Theme" theme name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (Name) translations: Spanish -> Portuguese -----------------
$ws.Range("C2").Value  = "Redstone Keep: Bazar dos Artesãos"
$ws.Range("C3").Value  = "Vale Frost: Bazar da Rosa Branca"
$ws.Range("C4").Value  = "Diresprings: Bazar do Deserto"
$ws.Range("C5").Value  = "Fleur: Bazar de Flores"
$ws.Range("C6").Value  = "Ponto de interesse"
$ws.Range("C7").Value  = "Comerciante de Habilidades"
$ws.Range("C8").Value  = "Comerciante Viajante 1"
$ws.Range("C9").Value  = "Comerciante Viajante 2"
$ws.Range("C10").Value = "Comerciante Viajante Sênior 1"
$ws.Range("C11").Value = "Comerciante de Bestas Comuns"
$ws.Range("C12").Value = "Comerciante de Bestas Comuns"
$ws.Range("C13").Value = "Mercenário Rosa Branca Comum"
$ws.Range("C14").Value = "Mercenário Mastigure Comum"
$ws.Range("C15").Value = "Mercenário Brasa Comum"

# --- Column C width: widen to fit the longer Portuguese strings ----------
$ws.Columns.Item(3).ColumnWidth = 33.28515625

# --- Selection: now spans the whole translated column --------------------
[void]$ws.Range("C2:C15").Select()

# --- Theme: restore the default Excel theme name -------------------------
$wb.Theme.Name = "Office Theme"

# --- Window geometry (best effort; matches the saved workbookView) -------
$win = $excel.ActiveWindow
$win.Left   = 1815
$win.Top    = 750
$win.Width  = 15375
$win.Height = 8325
